$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.782.97'
$ws.Range('E2').Value = '  -0.20%  '
$ws.Range('D3').Value = '2.327.36'
$ws.Range('E3').Value = '  +3.13%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '232.83'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +1.14%  '
$ws.Range('E6').Value = '  +1.33%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '66.24'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  +5.39%  '
$ws.Range('E8').Value = '  -0.07%  '
$ws.Range('E9').Value = '  -1.22%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0972'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -4.48%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '56.94'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -0.45%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '26.92'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  +4.14%  '
$ws.Range('D13').Value = '2.671.23'
$ws.Range('E13').Value = '  +3.06%  '
$ws.Range('E14').Value = '  -1.66%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '15.31'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  -1.99%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '6.13'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  -0.64%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.838'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  -1.12%  '
$ws.Range('D18').Value = '2.322.31'
$ws.Range('E18').Value = '  +3.44%  '
$ws.Range('D19').Value = '43.708.02'
$ws.Range('E19').Value = '  -0.12%  '
$ws.Range('D20').Value = '0.0₃0979'
$ws.Range('E20').Value = '  -3.06%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '73.55'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +0.24%  '
$ws.Range('E22').Value = '  +2.28%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '248.52'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -1.48%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '3.82'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +15.05%  '
$ws.Range('E25').Value = '  -0.13%  '
$ws.Range('E26').Value = '  +0.17%  '
$ws.Range('E27').Value = '  -2.85%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.86'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -1.60%  '
$ws.Range('B29').Value = 'EthereumClassic'
$ws.Range('C29').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '22.29'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +7.37%  '
$ws.Range('B30').Value = 'Monero'
$ws.Range('C30').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '174.81'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +1.70%  '
$ws.Range('B31').Value = 'ImmutableX'
$ws.Range('C31').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.43'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +3.62%  '
$ws.Range('B32').Value = 'Kaspa'
$ws.Range('C32').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.131'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -4.94%  '
$ws.Range('E34').Value = '  +4.20%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0684'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -0.54%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '4.98'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +2.39%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '6.54'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +0.65%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.67'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -4.14%  '
$ws.Range('E39').Value = '  +2.75%  '
$ws.Range('E40').Value = '  -1.63%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '9.12'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +10.47%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.999'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -0.18%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '17.38'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -0.82%  '
$ws.Range('B44').Value = 'Aave'
$ws.Range('C44').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '98.51'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +0.56%  '
$ws.Range('B45').Value = 'ARBITRUM'
$ws.Range('C45').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.14'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +6.31%  '
$ws.Range('E46').Value = '  +3.97%  '
$ws.Range('B47').Value = 'FTXToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '4.41'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +1.95%  '
$ws.Range('B48').Value = 'Cronos'
$ws.Range('C48').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0952'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -2.29%  '
$ws.Range('B49').Value = 'TrustWalletToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.19'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -0.97%  '
$ws.Range('E50').Value = '  -5.45%  '
$ws.Range('D51').Value = '1.440.38'
$ws.Range('E51').Value = '  -0.30%  '
